$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-71 down to 51-72
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new record
$ws.Cells.Item(50, 1).Value = 5
$ws.Cells.Item(50, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(50, 3).Value = "Maule"
$ws.Cells.Item(50, 4).Value = 44455
$ws.Cells.Item(50, 4).NumberFormat = $ws.Cells.Item(51, 4).NumberFormat
$ws.Cells.Item(50, 5).Value = 7
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100108
$ws.Cells.Item(50, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(50, 9).Value = 100108002
$ws.Cells.Item(50, 10).Value = "Mango"
$ws.Cells.Item(50, 11).Value = "Sin especificar"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 100
$ws.Cells.Item(50, 14).Value = 8000
$ws.Cells.Item(50, 15).Value = 8000
$ws.Cells.Item(50, 16).Value = 8000
$ws.Cells.Item(50, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(50, 18).Value = "Ecuador"
$ws.Cells.Item(50, 19).Value = 2000
$ws.Cells.Item(50, 20).Value = 4
